$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph "Action; All moves are made wit" + bookmark(_GoBack) + "h no
#    winner" -> merge the two runs into a single run with the full text and
#    drop the bookmark from this location (it moves to the list item below).
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(23)
$r1 = $p1.Range
$r1.End = $r1.End - 1                 # exclude the paragraph mark
$text1 = $r1.Text
$r1.Text = ""
$r1.InsertAfter($text1)
$r1.LanguageID = "en-CA"

# ---------------------------------------------------------------------------
# 2) Paragraph "Action; Test case for each of the 49 winning combinations."
#    gains a leading space, and the _GoBack bookmark is (re)created right
#    after that new space, splitting the text into two runs.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(26)
$insertPos = $p2.Range.Start

$sel = $word.Selection
$sel.SetRange($insertPos, $insertPos)
$sel.TypeText(" ")

$bmPos = $insertPos + 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 3) Paragraph "Finish implementing " / "Score " / "button " -> merge the
#    three runs into a single run.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(39)
$r3 = $p3.Range
$r3.End = $r3.End - 1
$text3 = $r3.Text
$r3.Text = ""
$r3.InsertAfter($text3)
$r3.LanguageID = "en-CA"

Write-Host "done"
